$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.036.61'
$ws.Range("E2").Value = '  +0.30%  '
$ws.Range("D3").Value = '2.534.68'
$ws.Range("E3").Value = '  -2.62%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '591.20'
$ws.Range("E5").Value = '  +1.33%  '
$ws.Range("D6").Value = '173.26'
$ws.Range("E6").Value = '  +4.63%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '0.528'
$ws.Range("E8").Value = '  +0.25%  '
$ws.Range("D9").Value = '2.537.41'
$ws.Range("E9").Value = '  -2.43%  '
$ws.Range("E10").Value = '  -1.21%  '
$ws.Range("E11").Value = '  +1.82%  '
$ws.Range("D12").Value = '5.14'
$ws.Range("E12").Value = '  -0.91%  '
$ws.Range("D13").Value = '0.344'
$ws.Range("E13").Value = '  -5.64%  '
$ws.Range("D14").Value = '26.73'
$ws.Range("E14").Value = '  -1.21%  '
$ws.Range("D15").Value = '3.009.76'
$ws.Range("E15").Value = '  -2.17%  '
$ws.Range("D16").Value = '0.0000176'
$ws.Range("E16").Value = '  -1.19%  '
$ws.Range("D17").Value = '66.910.15'
$ws.Range("E17").Value = '  -0.03%  '
$ws.Range("D18").Value = '2.559.36'
$ws.Range("E18").Value = '  -1.40%  '
$ws.Range("D19").Value = '8.05'
$ws.Range("E19").Value = '  +3.55%  '
$ws.Range("D20").Value = '11.30'
$ws.Range("E20").Value = '  -2.76%  '
$ws.Range("D21").Value = '354.02'
$ws.Range("E21").Value = '  +0.44%  '
$ws.Range("D22").Value = '4.18'
$ws.Range("E22").Value = '  -1.54%  '
$ws.Range("D23").Value = '4.60'
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").Value = '1.99'
$ws.Range("E24").Value = '  +4.85%  '
$ws.Range("D26").Value = '69.79'
$ws.Range("E26").Value = '  +1.09%  '
$ws.Range("D27").Value = '10.07'
$ws.Range("E27").Value = '  -3.39%  '
$ws.Range("D28").Value = '2.679.24'
$ws.Range("E28").Value = '  -2.09%  '
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("D30").Value = '0.0₃0984'
$ws.Range("E30").Value = '  -0.22%  '
$ws.Range("D31").Value = '533.18'
$ws.Range("E31").Value = '  -0.92%  '
$ws.Range("D32").Value = '8.14'
$ws.Range("E32").Value = '  +0.31%  '
$ws.Range("D33").Value = '1.34'
$ws.Range("E33").Value = '  +0.88%  '
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("D35").Value = '0.131'
$ws.Range("E35").Value = '  -0.97%  '
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("E37").Value = '  -0.25%  '
$ws.Range("D38").Value = '157.18'
$ws.Range("E38").Value = '  -0.45%  '
$ws.Range("D39").Value = '18.60'
$ws.Range("E39").Value = '  -1.14%  '
$ws.Range("D40").Value = '18.44'
$ws.Range("E40").Value = '  +1.04%  '
$ws.Range("D41").Value = '0.355'
$ws.Range("E41").Value = '  -1.84%  '
$ws.Range("E42").Value = '  -0.35%  '
$ws.Range("D43").Value = '5.12'
$ws.Range("E43").Value = '  +0.39%  '
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("D45").Value = '2.50'
$ws.Range("E45").Value = '  +4.27%  '
$ws.Range("D46").Value = '39.71'
$ws.Range("E46").Value = '  -1.11%  '
$ws.Range("D47").Value = '149.14'
$ws.Range("E47").Value = '  -0.57%  '
$ws.Range("D48").Value = '0.558'
$ws.Range("E48").Value = '  -2.45%  '
$ws.Range("D49").Value = '0.0₆0278'
$ws.Range("E49").Value = '  -4.71%  '
$ws.Range("D50").Value = '3.68'
$ws.Range("E50").Value = '  -1.45%  '
$ws.Range("D51").Value = '1.69'
$ws.Range("E51").Value = '  -0.37%  '
